$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Holidays 2019')

$ws.Range('C2').Value = '30 42 сер груз'
$ws.Range('G2').Value = '30, 42, сер, груз'
$ws.Range('C3').Value = 'сер легк б/к'
$ws.Range('G3').Value = '30, 42, сер, груз'
$ws.Range('C4').Value = 'сер легк б/к'
$ws.Range('G4').Value = 'сер, легк, б/к'
$ws.Range('C5').Value = 'сер ошип'
$ws.Range('G5').Value = 'сер, легк, б/к'
$ws.Range('C6').Value = 'сер легк'
$ws.Range('G6').Value = 'сер, ошип'
$ws.Range('C7').Value = 'сер легк'
$ws.Range('G7').Value = 'сер, легк'
$ws.Range('C8').Value = 'сер легк'
$ws.Range('G8').Value = 'сер, легк'
$ws.Range('C9').Value = 'сер легк'
$ws.Range('G9').Value = 'сер, легк'
$ws.Range('C10').Value = '210B сер Type H C'
$ws.Range('G10').Value = 'сер, легк'
$ws.Range('C11').Value = 'сер груз Type LS-2'
$ws.Range('G11').Value = '210B, сер, Type, H, C'
$ws.Range('C12').Value = '202B сер Type C'
$ws.Range('G12').Value = '210B, сер, Type, H, C'
$ws.Range('C13').Value = '202B сер Type LS-2 H C'
$ws.Range('G13').Value = 'сер, груз, Type, LS-2'
$ws.Range('C14').Value = 'сер груз б/к'
$ws.Range('G14').Value = '202B, сер, Type, C'
$ws.Range('C15').Value = 'сер легк б/к'
$ws.Range('G15').Value = '202B, сер, Type, LS-2, H, C'
$ws.Range('C16').Value = 'сер легк б/к'
$ws.Range('G16').Value = '202B, сер, Type, LS-2, H, C'
$ws.Range('G17').Value = '202B, сер, Type, LS-2, H, C'
$ws.Range('C18').Value = 'сер легк'
$ws.Range('G18').Value = 'сер, груз, б/к'
$ws.Range('C19').Value = 'сер легк'
$ws.Range('G19').Value = 'сер, груз, б/к'
$ws.Range('C20').Value = 'сер легк'
$ws.Range('G20').Value = 'сер, груз, б/к'
$ws.Range('C21').Value = 'сер легк'
$ws.Range('G21').Value = 'сер, груз, б/к'
$ws.Range('C22').Value = 'сер легк'
$ws.Range('G22').Value = 'сер, легк, б/к'
$ws.Range('G23').Value = 'сер, легк, б/к'
$ws.Range('G24').Value = 'сер, легк'
$ws.Range('G25').Value = 'сер, легк'
$ws.Range('G26').Value = 'сер, легк'
$ws.Range('G27').Value = 'сер, легк'
$ws.Range('G28').Value = 'сер, легк'
$ws.Range('G29').Value = 'сер, легк'
$ws.Range('G30').Value = 'сер, легк'
